$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recruitment name in A2 ("Recruiter Specialist" -> "ASE 2022")
$ws.Range("A2").Value = "ASE 2022"

# Adjust column A width as in the diff (18.6640625 -> 28)
$ws.Columns.Item(1).ColumnWidth = 27.166666666666668

# Move the active selection as recorded in the diff (C3 -> B12)
$ws.Range("B12").Select() | Out-Null
